{"js": "// Replace the date line and each two-digit multiplication problem with its\n// new value. Every \"before\" string below occurs exactly once in the\n// document, so a plain case-sensitive search + full-text replace is safe.\nconst replacements = [\n  [\"2025-08-12 Tuesday\", \"2025-08-13 Wednesday\"],\n  [\"76\u00d787=\", \"11\u00d763=\"],\n  [\"54\u00d756=\", \"73\u00d788=\"],\n  [\"35\u00d784=\", \"58\u00d741=\"],\n  [\"44\u00d773=\", \"52\u00d714=\"],\n  [\"25\u00d789=\", \"37\u00d737=\"],\n  [\"55\u00d712=\", \"17\u00d729=\"],\n  [\"43\u00d745=\", \"66\u00d768=\"],\n  [\"47\u00d796=\", \"69\u00d793=\"],\n  [\"80\u00d732=\", \"79\u00d792=\"],\n  [\"85\u00d798=\", \"58\u00d717=\"],\n  [\"26\u00d778=\", \"41\u00d744=\"],\n  [\"62\u00d712=\", \"91\u00d769=\"],\n  [\"79\u00d763=\", \"28\u00d793=\"],\n  [\"91\u00d712=\", \"40\u00d772=\"],\n  [\"37\u00d756=\", \"36\u00d754=\"],\n  [\"34\u00d735=\", \"75\u00d766=\"],\n  [\"29\u00d719=\", \"36\u00d762=\"],\n  [\"80\u00d715=\", \"98\u00d716=\"],\n  [\"14\u00d796=\", \"81\u00d754=\"],\n  [\"75\u00d728=\", \"20\u00d761=\"],\n  [\"32\u00d756=\", \"88\u00d794=\"],\n  [\"90\u00d713=\", \"67\u00d752=\"],\n  [\"32\u00d764=\", \"97\u00d712=\"],\n  [\"81\u00d799=\", \"38\u00d752=\"],\n  [\"34\u00d721=\", \"79\u00d778=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each two-digit multiplication problem with its\n# new value. Every \"before\" string below occurs exactly once in the\n# document, so Find/Replace (wdReplaceAll) scoped to the whole document is\n# safe and unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-08-12 Tuesday\", \"2025-08-13 Wednesday\"),\n  @(\"76\u00d787=\", \"11\u00d763=\"),\n  @(\"54\u00d756=\", \"73\u00d788=\"),\n  @(\"35\u00d784=\", \"58\u00d741=\"),\n  @(\"44\u00d773=\", \"52\u00d714=\"),\n  @(\"25\u00d789=\", \"37\u00d737=\"),\n  @(\"55\u00d712=\", \"17\u00d729=\"),\n  @(\"43\u00d745=\", \"66\u00d768=\"),\n  @(\"47\u00d796=\", \"69\u00d793=\"),\n  @(\"80\u00d732=\", \"79\u00d792=\"),\n  @(\"85\u00d798=\", \"58\u00d717=\"),\n  @(\"26\u00d778=\", \"41\u00d744=\"),\n  @(\"62\u00d712=\", \"91\u00d769=\"),\n  @(\"79\u00d763=\", \"28\u00d793=\"),\n  @(\"91\u00d712=\", \"40\u00d772=\"),\n  @(\"37\u00d756=\", \"36\u00d754=\"),\n  @(\"34\u00d735=\", \"75\u00d766=\"),\n  @(\"29\u00d719=\", \"36\u00d762=\"),\n  @(\"80\u00d715=\", \"98\u00d716=\"),\n  @(\"14\u00d796=\", \"81\u00d754=\"),\n  @(\"75\u00d728=\", \"20\u00d761=\"),\n  @(\"32\u00d756=\", \"88\u00d794=\"),\n  @(\"90\u00d713=\", \"67\u00d752=\"),\n  @(\"32\u00d764=\", \"97\u00d712=\"),\n  @(\"81\u00d799=\", \"38\u00d752=\"),\n  @(\"34\u00d721=\", \"79\u00d778=\")\n)\n\nforeach ($pair in $pairs) {\n  $find = $pair[0]\n  $replace = $pair[1]\n  $rng = $d.Content\n  $null = $rng.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
